$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4365677535533905
$ws.Range("B1").Value = 0.550317645072937
$ws.Range("C1").Value = 0.799268364906311
$ws.Range("D1").Value = 3.762670040130615
$ws.Range("E1").Value = 5.652971267700195
